$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 12 - "TextBox 5": reposition/resize so the wording still fits, and
# clarify the caption text.
# ---------------------------------------------------------------------------
$s12 = $p.Slides.Item(12)

$noteBox = $s12.Shapes.Item(5)
$noteBox.Left  = 235.56071
$noteBox.Width = 297.0227
$noteBox.TextFrame.TextRange.Text = "Note computation of string size."

# Slide 12 - "Connector: Elbow 8": flip it horizontally as well as
# vertically, and nudge its (essentially zero) width by a single EMU.
$connector = $s12.Shapes.Item(7)
$connector.HorizontalFlip = $true
$connector.Width = 0.0001

# ---------------------------------------------------------------------------
# Slide 14 - "TextBox 1": shift it to the left.
# ---------------------------------------------------------------------------
$s14 = $p.Slides.Item(14)
$s14.Shapes.Item(5).Left = 97.259

# ---------------------------------------------------------------------------
# Slide 2 - "Content Placeholder 2": tidy up capitalization/punctuation.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange
$tr2.Paragraphs(2, 1).Text = "String literals are enclosed in (double) quotation marks."
$tr2.Paragraphs(4, 1).Text = "A string variable has two integer properties."

# ---------------------------------------------------------------------------
# Slide 8 - "Content Placeholder 2": remove a stray extra space.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$tr8 = $s8.Shapes.Item(2).TextFrame.TextRange
$tr8.Paragraphs(4, 1).Runs(1, 1).Text = "variable = ( "
